# Add new vocabulary entries (Vocabs from Context and Meaning 2) into
# columns H (Vocab) and I (Meaning) for rows 17-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H17").Value = "Jeopardize "
$ws.Range("I17").Value = "risk damaging or destroying sth important"

$ws.Range("H18").Value = "hinder"
$ws.Range("I18").Value = "stop sth (from developing/progressing)"

$ws.Range("H19").Value = "disparate"
$ws.Range("I19").Value = "having many differences"

$ws.Range("H20").Value = "coarse"
$ws.Range("I20").Value = "rough and hard"

# Update the active cell selection to match the saved workbook state.
$ws.Range("I21").Select()
